# Automatische test-sync: 2025-08-28 18:05:50
#
# Append the new "Retour status" mail-log entry (row 7) to the Logs sheet,
# extend the conditional-formatting ranges so they keep covering the whole
# data range (D/G/H/I/J columns), and bump the Dashboard roll-up count for
# the "Retour / Terugbetaling" category from 5 to 6.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ---------------------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A7").Value = "Retour status"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("D7").Value = "Retour / Terugbetaling"
$ws.Range("F7").Value = "2025-08-28 18:05:29"
$ws.Range("G7").Value = "Ja"
$ws.Range("H7").Value = "Nee"
$ws.Range("I7").Value = "Nee"
$ws.Range("J7").Value = "Nee"

# --- Logs sheet: widen conditional formatting to include the new row ------
$ws.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D7"))
$ws.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))
$ws.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H7"))
$ws.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I7"))
$ws.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J7"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" tally --------------
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B2").Value = 6
